$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "03-00-00 / 교육 예산 관리" row entirely; everything below shifts up one row.
$ws.Rows.Item(18).Delete()

# The "No" / work-level codes in column A are plain literal text (no formulas
# anywhere in this sheet), so after the shift they need to be re-written to
# match the renumbered sequence.
$codesA = @{
    18 = "03-00-00"
    19 = "03-01-00"
    20 = "03-01-01"
    21 = "03-01-02"
    22 = "03-01-03"
    23 = "03-02-00"
    24 = "03-02-01"
    25 = "03-02-02"
    26 = "03-02-03"
    27 = "03-02-04"
    28 = "03-02-05"
    29 = "03-02-06"
    30 = "03-02-07"
    31 = "03-02-08"
    32 = "03-02-09"
    33 = "03-03-00"
    34 = "04-00-00"
    35 = "04-01-00"
    36 = "04-02-00"
    37 = "04-03-00"
    38 = "04-04-00"
    39 = "05-00-00"
    40 = "06-00-00"
    41 = "06-01-00"
    42 = "06-02-00"
    43 = "06-03-00"
    44 = "06-04-00"
    45 = "06-05-00"
    46 = "06-06-00"
    47 = "06-07-00"
    48 = "06-08-00"
    49 = "06-09-00"
}
foreach ($r in $codesA.Keys) {
    $ws.Cells.Item($r, 1).Value = $codesA[$r]
}

# Append the new last task row ("07-00-00 / 자동화 프로그램 개발").
$ws.Cells.Item(50, 1).Value = "07-00-00"
$ws.Cells.Item(50, 2).Value = "자동화 프로그램 개발"
$ws.Cells.Item(50, 3).Value = " "
$ws.Cells.Item(50, 4).Value = " "
$ws.Cells.Item(50, 5).Value = "진행"
$ws.Cells.Item(50, 6).Value = " "
$ws.Cells.Item(50, 7).Value = " "
$ws.Cells.Item(50, 8).Value = " "
$ws.Cells.Item(50, 9).Value = 1
